function Set-TextValue {
    param($cellRef, $val)
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue "D2" "228.62"
Set-TextValue "G2" "6"
Set-TextValue "D3" "22.59"
Set-TextValue "G3" "6"
Set-TextValue "D4" "5.257"
Set-TextValue "G4" "6"
Set-TextValue "D5" "0.05571"
Set-TextValue "G5" "6"
Set-TextValue "D6" "3.386"
Set-TextValue "G6" "6"
Set-TextValue "D7" "6.465"
Set-TextValue "G7" "6"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D8" "0.7844"
$ws.Range("E8").Value = "7MXTokenMX"
Set-TextValue "G8" "6"
$ws.Range("B9").Value = "FTXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue "D9" "1.046"
$ws.Range("E9").Value = "8FTXTokenFTT"
Set-TextValue "G9" "6"
Set-TextValue "D10" "0.1377"
Set-TextValue "G10" "6"
Set-TextValue "D11" "0.07323"
Set-TextValue "G11" "6"
Set-TextValue "D12" "0.03166"
Set-TextValue "G12" "6"
Set-TextValue "D13" "0.02976"
Set-TextValue "G13" "6"
Set-TextValue "D14" "0.09271"
Set-TextValue "G14" "6"
Set-TextValue "D15" "0.001667"
Set-TextValue "G15" "6"
Set-TextValue "D16" "3.254"
Set-TextValue "G16" "6"
Set-TextValue "G17" "6"
Set-TextValue "D18" "0.0005808"
Set-TextValue "G18" "6"
Set-TextValue "D19" "0.006228"
Set-TextValue "G19" "6"
Set-TextValue "D20" "0.005221"
Set-TextValue "G20" "6"
Set-TextValue "G21" "6"
Set-TextValue "G22" "6"
Set-TextValue "D23" "3.952"
Set-TextValue "G23" "6"
Set-TextValue "G24" "6"
Set-TextValue "G25" "6"
Set-TextValue "D26" "0.1051"
Set-TextValue "G26" "6"
Set-TextValue "D27" "0.0005007"
$ws.Range("E27").Value = "26UpBotsUBXT"
Set-TextValue "G27" "6"
Set-TextValue "G28" "6"
Set-TextValue "G29" "6"
Set-TextValue "G30" "6"
Set-TextValue "G31" "6"
Set-TextValue "G32" "6"
Set-TextValue "G33" "6"
Set-TextValue "G34" "6"
Set-TextValue "G35" "6"
Set-TextValue "G36" "6"
Set-TextValue "G37" "6"
Set-TextValue "G38" "6"
Set-TextValue "G39" "6"
Set-TextValue "D40" "0.03995"
Set-TextValue "G40" "6"
Set-TextValue "D41" "0.006999"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
Set-TextValue "G41" "6"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1036"
$ws.Range("E42").Value = "41BKEXTokenBKK"
Set-TextValue "G42" "6"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.003244"
$ws.Range("E43").Value = "42CEJICEJI"
Set-TextValue "G43" "6"
Set-TextValue "D44" "0.009106"
Set-TextValue "G44" "6"
Set-TextValue "D45" "0.00005447"
Set-TextValue "G45" "6"
Set-TextValue "G46" "6"
Set-TextValue "D47" "0.7863"
Set-TextValue "G47" "6"
Set-TextValue "D48" "0.04212"
Set-TextValue "G48" "6"
Set-TextValue "D49" "0.00002103"
Set-TextValue "G49" "6"
Set-TextValue "G50" "6"
Set-TextValue "G51" "6"
